$d = $word.ActiveDocument

$null = $d.Content.Find.Execute("36+24=60", $true, $false, $false, $false, $false, $true, 1, $false, "31+40=71", 2)
$null = $d.Content.Find.Execute("89-66=23", $true, $false, $false, $false, $false, $true, 1, $false, "17+68=85", 2)
$null = $d.Content.Find.Execute("16+18=34", $true, $false, $false, $false, $false, $true, 1, $false, "32-13=19", 2)
$null = $d.Content.Find.Execute("75-29=46", $true, $false, $false, $false, $false, $true, 1, $false, "78-55=23", 2)
$null = $d.Content.Find.Execute("64-13=51", $true, $false, $false, $false, $false, $true, 1, $false, "95-4=91", 2)
$null = $d.Content.Find.Execute("70-41=29", $true, $false, $false, $false, $false, $true, 1, $false, "24+2=26", 2)
$null = $d.Content.Find.Execute("52+34=86", $true, $false, $false, $false, $false, $true, 1, $false, "91-74=17", 2)
$null = $d.Content.Find.Execute("36+12=48", $true, $false, $false, $false, $false, $true, 1, $false, "11+3=14", 2)
$null = $d.Content.Find.Execute("36-2=34", $true, $false, $false, $false, $false, $true, 1, $false, "59-37=22", 2)
$null = $d.Content.Find.Execute("56+14=70", $true, $false, $false, $false, $false, $true, 1, $false, "81-26=55", 2)
$null = $d.Content.Find.Execute("16+13=29", $true, $false, $false, $false, $false, $true, 1, $false, "76-43=33", 2)
$null = $d.Content.Find.Execute("48+36=84", $true, $false, $false, $false, $false, $true, 1, $false, "60-52=8", 2)
$null = $d.Content.Find.Execute("16+61=77", $true, $false, $false, $false, $false, $true, 1, $false, "81+0=81", 2)
$null = $d.Content.Find.Execute("87+1=88", $true, $false, $false, $false, $false, $true, 1, $false, "32-23=9", 2)
$null = $d.Content.Find.Execute("15+41=56", $true, $false, $false, $false, $false, $true, 1, $false, "85-38=47", 2)
$null = $d.Content.Find.Execute("45+20=65", $true, $false, $false, $false, $false, $true, 1, $false, "28+29=57", 2)
$null = $d.Content.Find.Execute("78-18=60", $true, $false, $false, $false, $false, $true, 1, $false, "30-19=11", 2)
$null = $d.Content.Find.Execute("67-62=5", $true, $false, $false, $false, $false, $true, 1, $false, "86-69=17", 2)
$null = $d.Content.Find.Execute("49-43=6", $true, $false, $false, $false, $false, $true, 1, $false, "13+62=75", 2)
$null = $d.Content.Find.Execute("53+22=75", $true, $false, $false, $false, $false, $true, 1, $false, "23+6=29", 2)
$null = $d.Content.Find.Execute("51-16=35", $true, $false, $false, $false, $false, $true, 1, $false, "61+38=99", 2)
$null = $d.Content.Find.Execute("69-18=51", $true, $false, $false, $false, $false, $true, 1, $false, "29+49=78", 2)
$null = $d.Content.Find.Execute("71-54=17", $true, $false, $false, $false, $false, $true, 1, $false, "90-11=79", 2)
$null = $d.Content.Find.Execute("9+44=53", $true, $false, $false, $false, $false, $true, 1, $false, "70-16=54", 2)
$null = $d.Content.Find.Execute("4+30=34", $true, $false, $false, $false, $false, $true, 1, $false, "26+53=79", 2)
$null = $d.Content.Find.Execute("28+61=89", $true, $false, $false, $false, $false, $true, 1, $false, "80-10=70", 2)
$null = $d.Content.Find.Execute("89-82=7", $true, $false, $false, $false, $false, $true, 1, $false, "78-22=56", 2)
$null = $d.Content.Find.Execute("84+9=93", $true, $false, $false, $false, $false, $true, 1, $false, "47+38=85", 2)
$null = $d.Content.Find.Execute("59-46=13", $true, $false, $false, $false, $false, $true, 1, $false, "82-37=45", 2)
$null = $d.Content.Find.Execute("4+78=82", $true, $false, $false, $false, $false, $true, 1, $false, "16+26=42", 2)
$null = $d.Content.Find.Execute("6-5=1", $true, $false, $false, $false, $false, $true, 1, $false, "22+30=52", 2)
$null = $d.Content.Find.Execute("59+29=88", $true, $false, $false, $false, $false, $true, 1, $false, "27+51=78", 2)
$null = $d.Content.Find.Execute("25+72=97", $true, $false, $false, $false, $false, $true, 1, $false, "38+35=73", 2)
$null = $d.Content.Find.Execute("13+9=22", $true, $false, $false, $false, $false, $true, 1, $false, "0+27=27", 2)
$null = $d.Content.Find.Execute("22+55=77", $true, $false, $false, $false, $false, $true, 1, $false, "37+24=61", 2)
$null = $d.Content.Find.Execute("54-14=40", $true, $false, $false, $false, $false, $true, 1, $false, "31-18=13", 2)
$null = $d.Content.Find.Execute("8+50=58", $true, $false, $false, $false, $false, $true, 1, $false, "35+33=68", 2)
$null = $d.Content.Find.Execute("68-8=60", $true, $false, $false, $false, $false, $true, 1, $false, "91-78=13", 2)
$null = $d.Content.Find.Execute("30+63=93", $true, $false, $false, $false, $false, $true, 1, $false, "49-15=34", 2)
$null = $d.Content.Find.Execute("38-20=18", $true, $false, $false, $false, $false, $true, 1, $false, "53-0=53", 2)
$null = $d.Content.Find.Execute("62+4=66", $true, $false, $false, $false, $false, $true, 1, $false, "74-28=46", 2)
$null = $d.Content.Find.Execute("12+33=45", $true, $false, $false, $false, $false, $true, 1, $false, "52+16=68", 2)
$null = $d.Content.Find.Execute("78+17=95", $true, $false, $false, $false, $false, $true, 1, $false, "12+87=99", 2)
$null = $d.Content.Find.Execute("71-67=4", $true, $false, $false, $false, $false, $true, 1, $false, "20+60=80", 2)
$null = $d.Content.Find.Execute("71-19=52", $true, $false, $false, $false, $false, $true, 1, $false, "92-64=28", 2)
$null = $d.Content.Find.Execute("18+53=71", $true, $false, $false, $false, $false, $true, 1, $false, "62-9=53", 2)
$null = $d.Content.Find.Execute("68-65=3", $true, $false, $false, $false, $false, $true, 1, $false, "6+75=81", 2)
$null = $d.Content.Find.Execute("50+5=55", $true, $false, $false, $false, $false, $true, 1, $false, "76-71=5", 2)
$null = $d.Content.Find.Execute("90-68=22", $true, $false, $false, $false, $false, $true, 1, $false, "73+13=86", 2)
$null = $d.Content.Find.Execute("49+39=88", $true, $false, $false, $false, $false, $true, 1, $false, "72+19=91", 2)
$null = $d.Content.Find.Execute("1+88=89", $true, $false, $false, $false, $false, $true, 1, $false, "77-36=41", 2)
$null = $d.Content.Find.Execute("90-75=15", $true, $false, $false, $false, $false, $true, 1, $false, "81-21=60", 2)
$null = $d.Content.Find.Execute("78-4=74", $true, $false, $false, $false, $false, $true, 1, $false, "6+93=99", 2)
$null = $d.Content.Find.Execute("64-11=53", $true, $false, $false, $false, $false, $true, 1, $false, "46-8=38", 2)
$null = $d.Content.Find.Execute("72-63=9", $true, $false, $false, $false, $false, $true, 1, $false, "78-53=25", 2)
$null = $d.Content.Find.Execute("3+76=79", $true, $false, $false, $false, $false, $true, 1, $false, "33+48=81", 2)
$null = $d.Content.Find.Execute("62-53=9", $true, $false, $false, $false, $false, $true, 1, $false, "43+12=55", 2)
$null = $d.Content.Find.Execute("33+60=93", $true, $false, $false, $false, $false, $true, 1, $false, "14+63=77", 2)
$null = $d.Content.Find.Execute("67-18=49", $true, $false, $false, $false, $false, $true, 1, $false, "7-2=5", 2)
$null = $d.Content.Find.Execute("7+78=85", $true, $false, $false, $false, $false, $true, 1, $false, "51+47=98", 2)
$null = $d.Content.Find.Execute("21+10=31", $true, $false, $false, $false, $false, $true, 1, $false, "62+33=95", 2)
$null = $d.Content.Find.Execute("4+92=96", $true, $false, $false, $false, $false, $true, 1, $false, "57+26=83", 2)
$null = $d.Content.Find.Execute("46-4=42", $true, $false, $false, $false, $false, $true, 1, $false, "18+31=49", 2)
$null = $d.Content.Find.Execute("61+22=83", $true, $false, $false, $false, $false, $true, 1, $false, "56+12=68", 2)
$null = $d.Content.Find.Execute("9+37=46", $true, $false, $false, $false, $false, $true, 1, $false, "28-17=11", 2)
$null = $d.Content.Find.Execute("9+6=15", $true, $false, $false, $false, $false, $true, 1, $false, "87-37=50", 2)
$null = $d.Content.Find.Execute("21+6=27", $true, $false, $false, $false, $false, $true, 1, $false, "57-38=19", 2)
$null = $d.Content.Find.Execute("51+24=75", $true, $false, $false, $false, $false, $true, 1, $false, "4+44=48", 2)
$null = $d.Content.Find.Execute("21+23=44", $true, $false, $false, $false, $false, $true, 1, $false, "77-6=71", 2)
$null = $d.Content.Find.Execute("57-33=24", $true, $false, $false, $false, $false, $true, 1, $false, "55+39=94", 2)
$null = $d.Content.Find.Execute("75-27=48", $true, $false, $false, $false, $false, $true, 1, $false, "76+2=78", 2)
$null = $d.Content.Find.Execute("97-65=32", $true, $false, $false, $false, $false, $true, 1, $false, "27+7=34", 2)
$null = $d.Content.Find.Execute("4+60=64", $true, $false, $false, $false, $false, $true, 1, $false, "94-2=92", 2)
$null = $d.Content.Find.Execute("45+47=92", $true, $false, $false, $false, $false, $true, 1, $false, "99-84=15", 2)
$null = $d.Content.Find.Execute("77-7=70", $true, $false, $false, $false, $false, $true, 1, $false, "81-27=54", 2)
$null = $d.Content.Find.Execute("57-12=45", $true, $false, $false, $false, $false, $true, 1, $false, "11+67=78", 2)
$null = $d.Content.Find.Execute("20-15=5", $true, $false, $false, $false, $false, $true, 1, $false, "20+39=59", 2)
$null = $d.Content.Find.Execute("32+52=84", $true, $false, $false, $false, $false, $true, 1, $false, "93-68=25", 2)
$null = $d.Content.Find.Execute("14+1=15", $true, $false, $false, $false, $false, $true, 1, $false, "62-11=51", 2)
$null = $d.Content.Find.Execute("31-9=22", $true, $false, $false, $false, $false, $true, 1, $false, "37+39=76", 2)
$null = $d.Content.Find.Execute("84+8=92", $true, $false, $false, $false, $false, $true, 1, $false, "91-74=17", 2)
$null = $d.Content.Find.Execute("47-30=17", $true, $false, $false, $false, $false, $true, 1, $false, "60-17=43", 2)
$null = $d.Content.Find.Execute("34-34=0", $true, $false, $false, $false, $false, $true, 1, $false, "86-39=47", 2)
$null = $d.Content.Find.Execute("51+46=97", $true, $false, $false, $false, $false, $true, 1, $false, "25+40=65", 2)
$null = $d.Content.Find.Execute("22+20=42", $true, $false, $false, $false, $false, $true, 1, $false, "7-2=5", 2)
$null = $d.Content.Find.Execute("72-24=48", $true, $false, $false, $false, $false, $true, 1, $false, "39-2=37", 2)
$null = $d.Content.Find.Execute("94-66=28", $true, $false, $false, $false, $false, $true, 1, $false, "96-63=33", 2)
$null = $d.Content.Find.Execute("34-0=34", $true, $false, $false, $false, $false, $true, 1, $false, "87-1=86", 2)
$null = $d.Content.Find.Execute("19+26=45", $true, $false, $false, $false, $false, $true, 1, $false, "84-5=79", 2)
$null = $d.Content.Find.Execute("50+32=82", $true, $false, $false, $false, $false, $true, 1, $false, "57+37=94", 2)
$null = $d.Content.Find.Execute("77+9=86", $true, $false, $false, $false, $false, $true, 1, $false, "58+2=60", 2)
$null = $d.Content.Find.Execute("50+43=93", $true, $false, $false, $false, $false, $true, 1, $false, "72-29=43", 2)
$null = $d.Content.Find.Execute("49-47=2", $true, $false, $false, $false, $false, $true, 1, $false, "18+32=50", 2)
$null = $d.Content.Find.Execute("91-30=61", $true, $false, $false, $false, $false, $true, 1, $false, "92-68=24", 2)
$null = $d.Content.Find.Execute("47-24=23", $true, $false, $false, $false, $false, $true, 1, $false, "89-62=27", 2)
$null = $d.Content.Find.Execute("47+27=74", $true, $false, $false, $false, $false, $true, 1, $false, "68+20=88", 2)
$null = $d.Content.Find.Execute("67-39=28", $true, $false, $false, $false, $false, $true, 1, $false, "34+24=58", 2)
$null = $d.Content.Find.Execute("33+26=59", $true, $false, $false, $false, $false, $true, 1, $false, "77+4=81", 2)
$null = $d.Content.Find.Execute("55+6=61", $true, $false, $false, $false, $false, $true, 1, $false, "64-21=43", 2)
$null = $d.Content.Find.Execute("39-35=4", $true, $false, $false, $false, $false, $true, 1, $false, "89-31=58", 2)
